$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E13:H13 - new values for Melicertus kerathurus
$ws.Range("E13").Value = "Melicertus kerathurus"
$ws.Range("F13").Value = "MELIKER"
$ws.Range("G13").Value = 0.102
$ws.Range("H13").Value = 6

# Update E14:H14 - new values for Raja asterias
$ws.Range("E14").Value = "Raja asterias"
$ws.Range("F14").Value = "RAJAAST"
$ws.Range("G14").Value = 0.126
$ws.Range("H14").Value = 1

# Update E15:H15 - new values for Solea solea
$ws.Range("E15").Value = "Solea solea"
$ws.Range("F15").Value = "SOLEVUL"
$ws.Range("G15").Value = 2.488
$ws.Range("H15").Value = 20

# Update E16:H16 - new values for Squilla mantis; I16 also changes to 1
$ws.Range("E16").Value = "Squilla mantis"
$ws.Range("F16").Value = "SQUIMAN"
$ws.Range("G16").Value = 0.239
$ws.Range("H16").Value = 6
$ws.Range("I16").Value = 1

$wb.Save()
